# "Problem Solved Date input to Excel"
# The Timeline column (B2:B25) used to hold pre-computed date-label text
# (shared strings like "2015-01-31"). Replace each of those with a real
# =DATE(yyyy,mm,dd) formula so Excel treats the column as actual date
# input instead of static text. Also replace the placeholder row 26/27
# values with the real FORECAST.ETS() formula that predicts the next
# period from the Timeline/Room-Nights-Sold history, and drop the now
# unused row 27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(
    @(2015,1,31), @(2015,2,28), @(2015,3,31), @(2015,4,30),
    @(2015,5,31), @(2015,6,30), @(2015,7,31), @(2015,8,31),
    @(2015,9,30), @(2015,10,31), @(2015,11,30), @(2015,12,31),
    @(2016,1,31), @(2016,2,28), @(2016,3,31), @(2016,4,30),
    @(2016,5,31), @(2016,6,30), @(2016,7,31), @(2016,8,31),
    @(2016,9,30), @(2016,10,31), @(2016,11,30), @(2016,12,31)
)

$row = 2
foreach ($d in $dates) {
    $y = $d[0]
    $m = "{0:D2}" -f $d[1]
    $day = "{0:D2}" -f $d[2]
    $ws.Cells.Item($row, 2).Formula = "=DATE($y,$m,$day)"
    $row = $row + 1
}

# Row 26 (D/E/F used to each hold a literal 45 placeholder) becomes a
# single forecast formula in D26; E26/F26 are no longer used.
$ws.Cells.Item(26, 4).Formula = '=FORECAST.ETS(B26,$C$2:$C$25,$B$2:$B$25,1,1)'
$ws.Cells.Item(26, 5).ClearContents()
$ws.Cells.Item(26, 6).ClearContents()

# Row 27 (the duplicate D/E/F=45 placeholder row) is removed entirely.
$ws.Rows.Item(27).Delete()
